# Moving from 3.2.1 to 3.2.2: refresh the captured Java/JUnit/Tycho stack trace
# text inside the 'asTableNotExistingFile' expected-generation sample so it matches
# the exception dump produced by the newer toolchain (JDK module-qualified frames,
# renumbered JUnit/Surefire/Tycho/Equinox internals, Surefire/Tycho launch frames).

$d = $word.ActiveDocument

# The whole stack trace lives in a single run/Text element, so the safest way to
# apply the change is one literal whole-block Find & Replace (the block contains many
# repeated lines, e.g. ParentRunner$3.run(...), so line-by-line replace would be ambiguous).
$old = 'asTable(java.lang.String,java.lang.String,java.lang.String,java.lang.String,java.lang.String) with arguments [excel.xlsx, Feuil1, C3, F7, fr-FR] failed:{{NL}}{{TAB}}/home/development/git/M2Doc/tests/org.obeonetwork.m2doc.tests/resources/excelServices/asTableNotExistingFile/excel.xlsx (Aucun fichier ou dossier de ce type){{NL}}java.io.FileNotFoundException: /home/development/git/M2Doc/tests/org.obeonetwork.m2doc.tests/resources/excelServices/asTableNotExistingFile/excel.xlsx (Aucun fichier ou dossier de ce type){{NL}}{{TAB}}at java.io.FileInputStream.open0(Native Method){{NL}}{{TAB}}at java.io.FileInputStream.open(FileInputStream.java:195){{NL}}{{TAB}}at java.io.FileInputStream.<init>(FileInputStream.java:138){{NL}}{{TAB}}at org.eclipse.emf.ecore.resource.impl.FileURIHandlerImpl.createInputStream(FileURIHandlerImpl.java:99){{NL}}{{TAB}}at org.eclipse.emf.ecore.resource.impl.ExtensibleURIConverterImpl.createInputStream(ExtensibleURIConverterImpl.java:360){{NL}}{{TAB}}at org.eclipse.emf.ecore.resource.impl.ExtensibleURIConverterImpl.createInputStream(ExtensibleURIConverterImpl.java:354){{NL}}{{TAB}}at org.obeonetwork.m2doc.services.ExcelServices.asTable(ExcelServices.java:125){{NL}}{{TAB}}at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method){{NL}}{{TAB}}at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62){{NL}}{{TAB}}at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43){{NL}}{{TAB}}at java.lang.reflect.Method.invoke(Method.java:498){{NL}}{{TAB}}at org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:162){{NL}}{{TAB}}at org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:135){{NL}}{{TAB}}at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129){{NL}}{{TAB}}at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172){{NL}}{{TAB}}at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callOrApply(EvaluationServices.java:208){{NL}}{{TAB}}at org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:192){{NL}}{{TAB}}at org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119){{NL}}{{TAB}}at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53){{NL}}{{TAB}}at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69){{NL}}{{TAB}}at org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112){{NL}}{{TAB}}at org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:604){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:1){{NL}}{{TAB}}at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:186){{NL}}{{TAB}}at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53){{NL}}{{TAB}}at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1675){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1){{NL}}{{TAB}}at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199){{NL}}{{TAB}}at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53){{NL}}{{TAB}}at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:314){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1){{NL}}{{TAB}}at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279){{NL}}{{TAB}}at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53){{NL}}{{TAB}}at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:299){{NL}}{{TAB}}at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:853){{NL}}{{TAB}}at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:536){{NL}}{{TAB}}at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:422){{NL}}{{TAB}}at sun.reflect.GeneratedMethodAccessor6.invoke(Unknown Source){{NL}}{{TAB}}at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43){{NL}}{{TAB}}at java.lang.reflect.Method.invoke(Method.java:498){{NL}}{{TAB}}at org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:50){{NL}}{{TAB}}at org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12){{NL}}{{TAB}}at org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47){{NL}}{{TAB}}at org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17){{NL}}{{TAB}}at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27){{NL}}{{TAB}}at org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325){{NL}}{{TAB}}at org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78){{NL}}{{TAB}}at org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57){{NL}}{{TAB}}at org.junit.runners.ParentRunner$3.run(ParentRunner.java:290){{NL}}{{TAB}}at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71){{NL}}{{TAB}}at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288){{NL}}{{TAB}}at org.junit.runners.ParentRunner.access$000(ParentRunner.java:58){{NL}}{{TAB}}at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268){{NL}}{{TAB}}at org.junit.runners.ParentRunner.run(ParentRunner.java:363){{NL}}{{TAB}}at org.junit.runners.Suite.runChild(Suite.java:128){{NL}}{{TAB}}at org.junit.runners.Suite.runChild(Suite.java:27){{NL}}{{TAB}}at org.junit.runners.ParentRunner$3.run(ParentRunner.java:290){{NL}}{{TAB}}at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71){{NL}}{{TAB}}at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288){{NL}}{{TAB}}at org.junit.runners.ParentRunner.access$000(ParentRunner.java:58){{NL}}{{TAB}}at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268){{NL}}{{TAB}}at org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26){{NL}}{{TAB}}at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27){{NL}}{{TAB}}at org.junit.runners.ParentRunner.run(ParentRunner.java:363){{NL}}{{TAB}}at org.junit.runners.Suite.runChild(Suite.java:128){{NL}}{{TAB}}at org.junit.runners.Suite.runChild(Suite.java:27){{NL}}{{TAB}}at org.junit.runners.ParentRunner$3.run(ParentRunner.java:290){{NL}}{{TAB}}at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71){{NL}}{{TAB}}at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288){{NL}}{{TAB}}at org.junit.runners.ParentRunner.access$000(ParentRunner.java:58){{NL}}{{TAB}}at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268){{NL}}{{TAB}}at org.junit.runners.ParentRunner.run(ParentRunner.java:363){{NL}}{{TAB}}at org.junit.runners.Suite.runChild(Suite.java:128){{NL}}{{TAB}}at org.junit.runners.Suite.runChild(Suite.java:27){{NL}}{{TAB}}at org.junit.runners.ParentRunner$3.run(ParentRunner.java:290){{NL}}{{TAB}}at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71){{NL}}{{TAB}}at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288){{NL}}{{TAB}}at org.junit.runners.ParentRunner.access$000(ParentRunner.java:58){{NL}}{{TAB}}at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268){{NL}}{{TAB}}at org.junit.runners.ParentRunner.run(ParentRunner.java:363){{NL}}{{TAB}}at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86){{NL}}{{TAB}}at org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38){{NL}}{{TAB}}at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538){{NL}}{{TAB}}at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760){{NL}}{{TAB}}at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460){{NL}}{{TAB}}at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206){{NL}}'
$new = 'asTable(java.lang.String,java.lang.String,java.lang.String,java.lang.String,java.lang.String) with arguments [excel.xlsx, Feuil1, C3, F7, fr-FR] failed:{{NL}}{{TAB}}/home/development/git/M2Doc/tests/org.obeonetwork.m2doc.tests/resources/excelServices/asTableNotExistingFile/excel.xlsx (Aucun fichier ou dossier de ce type){{NL}}java.io.FileNotFoundException: /home/development/git/M2Doc/tests/org.obeonetwork.m2doc.tests/resources/excelServices/asTableNotExistingFile/excel.xlsx (Aucun fichier ou dossier de ce type){{NL}}{{TAB}}at java.base/java.io.FileInputStream.open0(Native Method){{NL}}{{TAB}}at java.base/java.io.FileInputStream.open(FileInputStream.java:212){{NL}}{{TAB}}at java.base/java.io.FileInputStream.<init>(FileInputStream.java:154){{NL}}{{TAB}}at org.eclipse.emf.ecore.resource.impl.FileURIHandlerImpl.createInputStream(FileURIHandlerImpl.java:99){{NL}}{{TAB}}at org.eclipse.emf.ecore.resource.impl.ExtensibleURIConverterImpl.createInputStream(ExtensibleURIConverterImpl.java:358){{NL}}{{TAB}}at org.eclipse.emf.ecore.resource.impl.ExtensibleURIConverterImpl.createInputStream(ExtensibleURIConverterImpl.java:352){{NL}}{{TAB}}at org.obeonetwork.m2doc.services.ExcelServices.asTable(ExcelServices.java:125){{NL}}{{TAB}}at java.base/jdk.internal.reflect.NativeMethodAccessorImpl.invoke0(Native Method){{NL}}{{TAB}}at java.base/jdk.internal.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62){{NL}}{{TAB}}at java.base/jdk.internal.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43){{NL}}{{TAB}}at java.base/java.lang.reflect.Method.invoke(Method.java:564){{NL}}{{TAB}}at org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:162){{NL}}{{TAB}}at org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:135){{NL}}{{TAB}}at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129){{NL}}{{TAB}}at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172){{NL}}{{TAB}}at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callOrApply(EvaluationServices.java:208){{NL}}{{TAB}}at org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189){{NL}}{{TAB}}at org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119){{NL}}{{TAB}}at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53){{NL}}{{TAB}}at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69){{NL}}{{TAB}}at org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:109){{NL}}{{TAB}}at org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:604){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:1){{NL}}{{TAB}}at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:186){{NL}}{{TAB}}at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53){{NL}}{{TAB}}at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1675){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1){{NL}}{{TAB}}at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199){{NL}}{{TAB}}at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53){{NL}}{{TAB}}at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:314){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1){{NL}}{{TAB}}at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279){{NL}}{{TAB}}at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53){{NL}}{{TAB}}at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450){{NL}}{{TAB}}at org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:299){{NL}}{{TAB}}at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:853){{NL}}{{TAB}}at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:529){{NL}}{{TAB}}at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:415){{NL}}{{TAB}}at jdk.internal.reflect.GeneratedMethodAccessor10.invoke(Unknown Source){{NL}}{{TAB}}at java.base/jdk.internal.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43){{NL}}{{TAB}}at java.base/java.lang.reflect.Method.invoke(Method.java:564){{NL}}{{TAB}}at org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:59){{NL}}{{TAB}}at org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12){{NL}}{{TAB}}at org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:56){{NL}}{{TAB}}at org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17){{NL}}{{TAB}}at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27){{NL}}{{TAB}}at org.junit.runners.ParentRunner$3.evaluate(ParentRunner.java:306){{NL}}{{TAB}}at org.junit.runners.BlockJUnit4ClassRunner$1.evaluate(BlockJUnit4ClassRunner.java:100){{NL}}{{TAB}}at org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:366){{NL}}{{TAB}}at org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:103){{NL}}{{TAB}}at org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:63){{NL}}{{TAB}}at org.junit.runners.ParentRunner$4.run(ParentRunner.java:331){{NL}}{{TAB}}at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:79){{NL}}{{TAB}}at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329){{NL}}{{TAB}}at org.junit.runners.ParentRunner.access$100(ParentRunner.java:66){{NL}}{{TAB}}at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:293){{NL}}{{TAB}}at org.junit.runners.ParentRunner.run(ParentRunner.java:413){{NL}}{{TAB}}at org.junit.runners.Suite.runChild(Suite.java:128){{NL}}{{TAB}}at org.junit.runners.Suite.runChild(Suite.java:27){{NL}}{{TAB}}at org.junit.runners.ParentRunner$4.run(ParentRunner.java:331){{NL}}{{TAB}}at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:79){{NL}}{{TAB}}at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329){{NL}}{{TAB}}at org.junit.runners.ParentRunner.access$100(ParentRunner.java:66){{NL}}{{TAB}}at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:293){{NL}}{{TAB}}at org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26){{NL}}{{TAB}}at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27){{NL}}{{TAB}}at org.junit.runners.ParentRunner$3.evaluate(ParentRunner.java:306){{NL}}{{TAB}}at org.junit.runners.ParentRunner.run(ParentRunner.java:413){{NL}}{{TAB}}at org.junit.runners.Suite.runChild(Suite.java:128){{NL}}{{TAB}}at org.junit.runners.Suite.runChild(Suite.java:27){{NL}}{{TAB}}at org.junit.runners.ParentRunner$4.run(ParentRunner.java:331){{NL}}{{TAB}}at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:79){{NL}}{{TAB}}at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329){{NL}}{{TAB}}at org.junit.runners.ParentRunner.access$100(ParentRunner.java:66){{NL}}{{TAB}}at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:293){{NL}}{{TAB}}at org.junit.runners.ParentRunner$3.evaluate(ParentRunner.java:306){{NL}}{{TAB}}at org.junit.runners.ParentRunner.run(ParentRunner.java:413){{NL}}{{TAB}}at org.junit.runners.Suite.runChild(Suite.java:128){{NL}}{{TAB}}at org.junit.runners.Suite.runChild(Suite.java:27){{NL}}{{TAB}}at org.junit.runners.ParentRunner$4.run(ParentRunner.java:331){{NL}}{{TAB}}at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:79){{NL}}{{TAB}}at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329){{NL}}{{TAB}}at org.junit.runners.ParentRunner.access$100(ParentRunner.java:66){{NL}}{{TAB}}at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:293){{NL}}{{TAB}}at org.junit.runners.ParentRunner$3.evaluate(ParentRunner.java:306){{NL}}{{TAB}}at org.junit.runners.ParentRunner.run(ParentRunner.java:413){{NL}}{{TAB}}at org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:365){{NL}}{{TAB}}at org.apache.maven.surefire.junit4.JUnit4Provider.executeWithRerun(JUnit4Provider.java:273){{NL}}{{TAB}}at org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:238){{NL}}{{TAB}}at org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:159){{NL}}{{TAB}}at java.base/jdk.internal.reflect.NativeMethodAccessorImpl.invoke0(Native Method){{NL}}{{TAB}}at java.base/jdk.internal.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62){{NL}}{{TAB}}at java.base/jdk.internal.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43){{NL}}{{TAB}}at java.base/java.lang.reflect.Method.invoke(Method.java:564){{NL}}{{TAB}}at org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:206){{NL}}{{TAB}}at org.apache.maven.surefire.booter.ProviderFactory$ProviderProxy.invoke(ProviderFactory.java:161){{NL}}{{TAB}}at org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:84){{NL}}{{TAB}}at org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:113){{NL}}{{TAB}}at org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21){{NL}}{{TAB}}at java.base/jdk.internal.reflect.NativeMethodAccessorImpl.invoke0(Native Method){{NL}}{{TAB}}at java.base/jdk.internal.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62){{NL}}{{TAB}}at java.base/jdk.internal.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43){{NL}}{{TAB}}at java.base/java.lang.reflect.Method.invoke(Method.java:564){{NL}}{{TAB}}at org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:593){{NL}}{{TAB}}at org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:205){{NL}}{{TAB}}at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:137){{NL}}{{TAB}}at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:107){{NL}}{{TAB}}at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:401){{NL}}{{TAB}}at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:255){{NL}}{{TAB}}at java.base/jdk.internal.reflect.NativeMethodAccessorImpl.invoke0(Native Method){{NL}}{{TAB}}at java.base/jdk.internal.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62){{NL}}{{TAB}}at java.base/jdk.internal.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43){{NL}}{{TAB}}at java.base/java.lang.reflect.Method.invoke(Method.java:564){{NL}}{{TAB}}at org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:657){{NL}}{{TAB}}at org.eclipse.equinox.launcher.Main.basicRun(Main.java:594){{NL}}{{TAB}}at org.eclipse.equinox.launcher.Main.run(Main.java:1447){{NL}}{{TAB}}at org.eclipse.equinox.launcher.Main.main(Main.java:1420){{NL}}'
$old = $old.Replace('{{TAB}}', [char]9).Replace('{{NL}}', [char]10)
$new = $new.Replace('{{TAB}}', [char]9).Replace('{{NL}}', [char]10)

$range = $d.Content
$replaced = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
Write-Output "Stack trace block replaced: $replaced"
